$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily record for 2021-02-08 (Excel serial 44235) was inserted into the
# table, which pushes the existing rows 93-113 down to 94-114, and a further
# new record for 2021-03-02 (serial 44257) was appended as the new last row
# (115). Because columns C/D are a trailing 7-entry rolling sum of column B
# (not live formulas - they're stored as plain numbers), every row from 90
# through the new last data row has to carry a refreshed total too.

# Row 114 and 115 do not exist yet in the sheet, so give them the same
# look (border / bold-ish date font / date number format / alignment) as the
# rest of the date column before writing values into them.
$ws.Range("A113:D113").Copy()
$ws.Range("A114:D115").PasteSpecial(-4122)

# Final contents for A90:D115 (dates as serials, "nuovi pos." counts, and the
# refreshed 7-day rolling totals / per-100k figures).
$data = @(
    @(90,  44232, 6,  24, 258.5092632485997),
    @(91,  44233, 6,  23, 247.7380439465747),
    @(92,  44234, 5,  21, 226.1956053425248),
    @(93,  44235, 4,  21, 226.1956053425248),
    @(94,  44236, 0,  18, 193.8819474364498),
    @(95,  44237, 0,  17, 183.1107281344248),
    @(96,  44238, 0,  21, 226.1956053425248),
    @(97,  44239, 3,  28, 301.5941404566997),
    @(98,  44240, 5,  30, 323.1365790607497),
    @(99,  44241, 9,  30, 323.1365790607497),
    @(100, 44242, 11, 37, 398.5351141749246),
    @(101, 44243, 2,  44, 473.9336492890995),
    @(102, 44244, 0,  58, 624.7307195174493),
    @(103, 44245, 7,  52, 560.1034037052995),
    @(104, 44246, 10, 56, 603.1882809133994),
    @(105, 44247, 19, 60, 646.2731581214994),
    @(106, 44248, 3,  64, 689.3580353295993),
    @(107, 44249, 15, 71, 764.7565704437743),
    @(108, 44250, 6,  70, 753.9853511417492),
    @(109, 44251, 4,  69, 743.2141318397242),
    @(110, 44252, 14, 72, 775.5277897457992),
    @(111, 44253, 9,  81, 872.4687634640242),
    @(112, 44254, 18, 80, 861.6975441619991),
    @(113, 44255, 6,  $null, $null),
    @(114, 44256, 24, $null, $null),
    @(115, 44257, 5,  $null, $null)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    if ($row[3] -eq $null) {
        $ws.Cells.Item($r, 3).ClearContents()
        $ws.Cells.Item($r, 4).ClearContents()
    } else {
        $ws.Cells.Item($r, 3).Value = $row[3]
        $ws.Cells.Item($r, 4).Value = $row[4]
    }
}
